# SCRC 2024 Olympics Pool PICKS - "Update - After 100m (Male)"
#
# After the Men's 100m Final was actually run, the picks that had been
# copy/pasted into the ALL_ENTRIES sheet (as a second "4x100m" block in
# rows 60-69) are removed, and Shawna's still-to-be-decided picks for the
# remaining events (100m - F, 200m - M/F, 400m - M/F, 4x100m - M/F) are
# cleared out on her own tab so she can fill them back in.

$wb = $excel.ActiveWorkbook

$allEntries = $wb.Worksheets.Item("ALL_ENTRIES")
$shawna     = $wb.Worksheets.Item("Shawna")

# ---------------------------------------------------------------------
# ALL_ENTRIES: the accidental duplicate header + 8 rows (rows 60-69) are
# deleted entirely, shrinking the sheet from A1:E69 down to A1:E59.
# ---------------------------------------------------------------------
$allEntries.Rows("60:69").Delete()

# Column E narrows slightly (was sized to fit the now-removed content).
$allEntries.Columns("E").ColumnWidth = 5.17

# ---------------------------------------------------------------------
# Shawna: her picks for B2:D9 (winner / 2nd / 3rd for each remaining
# event) are cleared out, leaving the event names in column A and her
# "locked-in" 100m pick notes in column E untouched.
# ---------------------------------------------------------------------
$shawna.Range("B2:D9").ClearContents()

# The columns re-fit themselves to the now-shorter contents.
$shawna.Columns("A").ColumnWidth = 35.67
$shawna.Columns("B").ColumnWidth = 12.0
$shawna.Columns("C").ColumnWidth = 15.5
$shawna.Columns("D").ColumnWidth = 14.83
$shawna.Columns("E").ColumnWidth = 12.83

# ---------------------------------------------------------------------
# Final selection / active sheet state as left by the editor.
# ---------------------------------------------------------------------
$allEntries.Range("B2:D9").Select() | Out-Null

$shawna.Activate() | Out-Null
$shawna.Range("C9").Select() | Out-Null
